$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Add six new rows (9-14) that mirror rows 3-8 (same Typ/Strategia/
#    Wskazniki/Warunki columns B:F) but with a new "1.x" numbering in
#    column A (new simulation series "3.1 - Simulation for complex
#    strategies, with Stop Loss").
# ---------------------------------------------------------------------------

# Row 9 <- copy of row 3, new label "1.1"
$ws.Range("B3:F3").Copy()
$ws.Range("B9:F9").PasteSpecial(-4122)
$ws.Range("B3:F3").Copy()
$ws.Range("B9:F9").PasteSpecial(-4163)
$ws.Range("A3").Copy()
$ws.Range("A9").PasteSpecial(-4122)
$ws.Range("A9").NumberFormat = "@"
$ws.Range("A9").Value = "1.1"

# Row 10 <- copy of row 4, new label "1.2"
$ws.Range("B4:F4").Copy()
$ws.Range("B10:F10").PasteSpecial(-4122)
$ws.Range("B4:F4").Copy()
$ws.Range("B10:F10").PasteSpecial(-4163)
$ws.Range("A8").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "1.2"

# Row 11 <- copy of row 5, new label "1.3"
$ws.Range("B5:F5").Copy()
$ws.Range("B11:F11").PasteSpecial(-4122)
$ws.Range("B5:F5").Copy()
$ws.Range("B11:F11").PasteSpecial(-4163)
$ws.Range("A3").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("A11").NumberFormat = "@"
$ws.Range("A11").Value = "1.3"

# Row 12 <- copy of row 6, new label "1.4"
$ws.Range("B6:F6").Copy()
$ws.Range("B12:F12").PasteSpecial(-4122)
$ws.Range("B6:F6").Copy()
$ws.Range("B12:F12").PasteSpecial(-4163)
$ws.Range("A8").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "1.4"

# Row 13 <- copy of row 7, new label "1.5"
$ws.Range("B7:F7").Copy()
$ws.Range("B13:F13").PasteSpecial(-4122)
$ws.Range("B7:F7").Copy()
$ws.Range("B13:F13").PasteSpecial(-4163)
$ws.Range("A3").Copy()
$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").NumberFormat = "@"
$ws.Range("A13").Value = "1.5"

# Row 14 <- copy of row 8, new label "1.6"
$ws.Range("B8:F8").Copy()
$ws.Range("B14:F14").PasteSpecial(-4122)
$ws.Range("B8:F8").Copy()
$ws.Range("B14:F14").PasteSpecial(-4163)
$ws.Range("A8").Copy()
$ws.Range("A14").PasteSpecial(-4122)
$ws.Range("A14").NumberFormat = "@"
$ws.Range("A14").Value = "1.6"

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2. Row heights - shrink the long-wrapped rows now that the strategy
#    descriptions have been condensed/re-flowed for the expanded table.
# ---------------------------------------------------------------------------
$ws.Rows("3:3").RowHeight = 135
$ws.Rows("4:4").RowHeight = 60
$ws.Rows("5:5").RowHeight = 45
$ws.Rows("6:6").RowHeight = 60
$ws.Rows("7:7").RowHeight = 45
$ws.Rows("8:8").RowHeight = 150
$ws.Rows("9:9").RowHeight = 135
$ws.Rows("10:10").RowHeight = 60
$ws.Rows("11:11").RowHeight = 45
$ws.Rows("12:12").RowHeight = 60
$ws.Rows("13:13").RowHeight = 45
$ws.Rows("14:14").RowHeight = 150

# ---------------------------------------------------------------------------
# 3. Page setup (print as A4 / portrait, as configured for the new table).
# ---------------------------------------------------------------------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# 4. View state - scroll down to show the new rows, select the last cell.
# ---------------------------------------------------------------------------
$ws.Range("F14").Select()
$excel.ActiveWindow.ScrollRow = 9
